$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "language" column header
$ws.Range("E1").Value = "language"

# Refresh randomized email addresses (column C) for rows 2-4
$ws.Range("C2").Value = "taotylwao@emltmp.com"
$ws.Range("C3").Value = "pshcamwao@emltmp.com"
$ws.Range("C4").Value = "ejkgnwcsn@emlpro.com"

# Update interest (column D) for rows 5-6
$ws.Range("D5").Value = "тварини"
$ws.Range("D6").Value = "hamster"

# Populate new language values (column E) for rows 2-6
$ws.Range("E2").Value = "en"
$ws.Range("E3").Value = "en"
$ws.Range("E4").Value = "en"
$ws.Range("E5").Value = "ua"
$ws.Range("E6").Value = "de"

# Row heights settle after the content/width change (auto-fit side effect)
$ws.Rows.Item(4).RowHeight = 19.5
$ws.Rows.Item(5).RowHeight = 20.25
